$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header style from an
# existing header cell (B1) so it reuses style index "1" (bold, centered,
# bordered) rather than creating a brand-new style.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F52).
$ws.Range("F2").Value = "2021-10-05 13:40:58.936018"
$ws.Range("F3").Value = "2021-10-05 13:40:58.936031"
$ws.Range("F4").Value = "2021-10-05 13:40:58.936035"
$ws.Range("F5").Value = "2021-10-05 13:40:58.936038"
$ws.Range("F6").Value = "2021-10-05 13:40:58.936040"
$ws.Range("F7").Value = "2021-10-05 13:40:58.936043"
$ws.Range("F8").Value = "2021-10-05 13:40:58.936046"
$ws.Range("F9").Value = "2021-10-05 13:40:58.936049"
$ws.Range("F10").Value = "2021-10-05 13:40:58.936052"
$ws.Range("F11").Value = "2021-10-05 13:40:58.936054"
$ws.Range("F12").Value = "2021-10-05 13:40:58.936057"
$ws.Range("F13").Value = "2021-10-05 13:40:58.936060"
$ws.Range("F14").Value = "2021-10-05 13:40:58.936062"
$ws.Range("F15").Value = "2021-10-05 13:40:58.936065"
$ws.Range("F16").Value = "2021-10-05 13:40:58.936067"
$ws.Range("F17").Value = "2021-10-05 13:40:58.936070"
$ws.Range("F18").Value = "2021-10-05 13:40:58.936073"
$ws.Range("F19").Value = "2021-10-05 13:40:58.936076"
$ws.Range("F20").Value = "2021-10-05 13:40:58.936078"
$ws.Range("F21").Value = "2021-10-05 13:40:58.936081"
$ws.Range("F22").Value = "2021-10-05 13:40:58.936084"
$ws.Range("F23").Value = "2021-10-05 13:40:58.936086"
$ws.Range("F24").Value = "2021-10-05 13:40:58.936089"
$ws.Range("F25").Value = "2021-10-05 13:40:58.936091"
$ws.Range("F26").Value = "2021-10-05 13:40:58.936094"
$ws.Range("F27").Value = "2021-10-05 13:40:58.936097"
$ws.Range("F28").Value = "2021-10-05 13:40:58.936099"
$ws.Range("F29").Value = "2021-10-05 13:40:58.936102"
$ws.Range("F30").Value = "2021-10-05 13:40:58.936104"
$ws.Range("F31").Value = "2021-10-05 13:40:58.936107"
$ws.Range("F32").Value = "2021-10-05 13:40:58.936109"
$ws.Range("F33").Value = "2021-10-05 13:40:58.936112"
$ws.Range("F34").Value = "2021-10-05 13:40:58.936115"
$ws.Range("F35").Value = "2021-10-05 13:40:58.936118"
$ws.Range("F36").Value = "2021-10-05 13:40:58.936120"
$ws.Range("F37").Value = "2021-10-05 13:40:58.936123"
$ws.Range("F38").Value = "2021-10-05 13:40:58.936126"
$ws.Range("F39").Value = "2021-10-05 13:40:58.936128"
$ws.Range("F40").Value = "2021-10-05 13:40:58.936131"
$ws.Range("F41").Value = "2021-10-05 13:40:58.936133"
$ws.Range("F42").Value = "2021-10-05 13:40:58.936136"
$ws.Range("F43").Value = "2021-10-05 13:40:58.936139"
$ws.Range("F44").Value = "2021-10-05 13:40:58.936142"
$ws.Range("F45").Value = "2021-10-05 13:40:58.936144"
$ws.Range("F46").Value = "2021-10-05 13:40:58.936147"
$ws.Range("F47").Value = "2021-10-05 13:40:58.936150"
$ws.Range("F48").Value = "2021-10-05 13:40:58.936153"
$ws.Range("F49").Value = "2021-10-05 13:40:58.936155"
$ws.Range("F50").Value = "2021-10-05 13:40:58.936158"
$ws.Range("F51").Value = "2021-10-05 13:40:58.936160"
$ws.Range("F52").Value = "2021-10-05 13:40:58.936163"
